# Yang Xing's week #1 Journal — apply the diff via Word COM interop.
#
# Strategy: each bullet that needs to grow a "detail" paragraph (or have its
# own paragraph properties/runs rewritten) is rewritten in place by feeding a
# small, well-formed OOXML package fragment to Range.InsertXML on that
# paragraph's Range. InsertXML replaces the targeted range with exactly the
# <w:p> elements supplied, so one call can both edit the existing paragraph
# and splice in brand-new sibling paragraphs right after it — precisely what
# the diff needs (numPr/pStyle kept intact, multi-run splits, proofErr tags,
# xml:space="preserve", custom indents, bookmarks, etc. all come through
# byte-for-byte). We walk the document from the LAST paragraph to the FIRST
# so that earlier paragraph indices never shift underneath us.

$d = $word.ActiveDocument

function New-BodyXml([string]$fragment) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $fragment + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function Set-ParagraphXml($paragraph, [string]$fragment) {
    $paragraph.Range.InsertXML((New-BodyXml $fragment))
}

# --- Paragraph 9 (trailing empty <w:p/>): delete it -------------------------
$d.Paragraphs(9).Range.Delete()

# --- Paragraph 8 "Study on the Node.js" -------------------------------------
# Drop the bookmark from this paragraph and give it its own following
# paragraph that holds just the bookmark (matches the diff's paragraph split).
$frag8 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Study on the Node.js</w:t></w:r></w:p>' +
    '<w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
Set-ParagraphXml $d.Paragraphs(8) $frag8

# --- Paragraph 7 "Study on the MongoDB" -------------------------------------
# Unchanged itself, but gains a new non-numbered detail paragraph after it.
$frag7 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Study on the MongoDB</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="705"/></w:pPr><w:r><w:t xml:space="preserve">Still need more test on MongoDB to make sure that it could meet our </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>requirements.</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>'
Set-ParagraphXml $d.Paragraphs(7) $frag7

# --- Paragraph 5 "Finish the Database schema design version 1..." ----------
# Becomes: "Make the draft..." bullet + its detail, then the original
# "Finish the Database..." bullet (with an added trailing ":" run) + its
# detail paragraph (with spell-check proofErr tags around the name).
$frag5 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Make the draft of the requirement document</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t>I wrote design document with some other group members. Please refer to the requirement document for the detail info.</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Finish the Database schema design version 1, write them down in the document</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t xml:space="preserve">I wrote the Database design document with </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Yuanyuan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Jia</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. We firstly list all the data Field according to our business logic without think whether we use SQL or NoSQL database to accomplish them. After that, we analysis the MongoDB technical characters and intend to use MongoDB as the project Database. Need to make further tests. </w:t></w:r></w:p>'
Set-ParagraphXml $d.Paragraphs(5) $frag5

# --- Paragraph 4 "Make the draft of the requirement document" --------------
# Becomes the "Design and settle down..." bullet (text + ":" run) plus its
# detail paragraph.
$frag4 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Design and settle down the business logic</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t xml:space="preserve">The design is under the assumption that there will be no more than one restaurant in one zip code scope. We mainly accomplish the menu display, dishes order, recommend dishes for customers, </w:t></w:r><w:r><w:t xml:space="preserve">make payment and </w:t></w:r><w:r><w:t xml:space="preserve">allow customers to make comments on the dishes they ordered. </w:t></w:r></w:p>'
Set-ParagraphXml $d.Paragraphs(4) $frag4

# --- Paragraph 3 "Design and settle down the business logic" ---------------
# Becomes "Fix down the project theme:" plus its detail paragraph.
$frag3 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Fix down the project theme:</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t xml:space="preserve">We decide to make an online dishes order restaurant which has multiple sites in different location.  </w:t></w:r></w:p>'
Set-ParagraphXml $d.Paragraphs(3) $frag3

Write-Output "Applied journal edits; final paragraph count:"
Write-Output $d.Paragraphs.Count
